# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds data rows 2-21 (20 rows of x,y,z gyroscope samples)
# beneath a header row. Two new samples are being inserted at the very top of
# the data (new rows 2-3), pushing the existing 20 rows down to rows 4-23, and
# eight new samples are appended after that (new rows 24-31).

# Shift the existing data rows (2-21) down by two rows, working from the
# bottom up so we never overwrite a row before it has been copied.
for ($r = 21; $r -ge 2; $r--) {
    $destRow = $r + 2
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value2
}

# Write the two new samples into the now-vacant rows 2 and 3.
$newTopRows = @(
    @(0.00006657902849829999, 0.5838314890861511, 0.3130545914173126),
    @(1.05201518535614, 1.37851881980896, 0.5420864224433899)
)
for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $newTopRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTopRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTopRows[$i][2]
}

# Append the eight new samples after the shifted data, in rows 24-31.
$newBottomRows = @(
    @(8.482767105102539, -0.4398876428604126, -8.189353942871094),
    @(-0.3088601231575012, -7.590475082397461, 2.290185451507568),
    @(6.10509729385376, -0.359460175037384, -4.315519332885742),
    @(-1.245227575302124, -1.861483097076416, 1.978062987327576),
    @(-1.619668006896973, 1.064798355102539, 1.856090188026428),
    @(-10.51302814483643, -17.02712059020996, -4.28569221496582),
    @(9.720071792602541, -2.732869386672974, -1.917076587677002),
    @(0.3920839130878448, -4.407997608184815, 1.71547520160675)
)
for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $r = 24 + $i
    $ws.Cells.Item($r, 1).Value = $newBottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottomRows[$i][2]
}
